# Add two new columns ("ligand_conc" and "time_step") to the studies endpoint
# template, inserted immediately before the existing "Comment" column (P),
# pushing Comment/Date/CTD Project Number/Collaborator two columns to the
# right (P:S -> R:U).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns at P:Q - this shifts the old P:S (Comment, Date,
# CTD Project Number, Collaborator) to R:U, carrying their values/styles
# along automatically (matches the diff for every data row 2-23 as well as
# the footer rows 26-27).
$ws.Range("P1:Q1").EntireColumn.Insert()

# New header labels for the freshly inserted columns.
$ws.Range("P1").Value = "ligand_conc"
$ws.Range("Q1").Value = "time_step"

# Give the two new columns the same width as the other general-purpose
# data columns (e.g. F:H), matching the template's styling.
$ws.Range("P1:Q1").ColumnWidth = 16.1953125

# Restore the selection to the (new) first of the shifted-right columns.
$ws.Range("Q1").Select()
